$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 9

$ws.Cells.Item($row, 1).Value = 1
$ws.Cells.Item($row, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item($row, 3).Value = "Arica y Parinacota"
$ws.Cells.Item($row, 4).Value = 44890
$ws.Cells.Item($row, 5).Value = 15
$ws.Cells.Item($row, 6).Value = "Fruta"
$ws.Cells.Item($row, 7).Value = 100108
$ws.Cells.Item($row, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item($row, 9).Value = 100108007
$ws.Cells.Item($row, 10).Value = "Coco"
$ws.Cells.Item($row, 11).Value = "Sin especificar"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 80
$ws.Cells.Item($row, 14).Value = 20000
$ws.Cells.Item($row, 15).Value = 23000
$ws.Cells.Item($row, 16).Value = 22250
$ws.Cells.Item($row, 17).Value = '$/malla 20 unidades'
$ws.Cells.Item($row, 18).Value = "Perú"
$ws.Cells.Item($row, 19).Value = 1112
$ws.Cells.Item($row, 20).Value = 20

$ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item($row - 1, 4).NumberFormat
